$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The task that was "EN PROCESO" (column B, row 36) has been finished and
# moves to "TERMINADAS" (column C, row 36) - keeping the same text/shared
# string, just relocated to the other column.
$finishedTask = $ws.Range("B36").Text
$ws.Range("C36").Value = $finishedTask
$ws.Range("B36").ClearContents()

# A new task enters "EN PROCESO" in row 37 (its cell already carries the
# correct style, so we only need to set its value).
$ws.Range("B37").Value = "Ajustar php word con los nuevos parametros de incrementos automaticos"

# Row 38 (now the next free "EN PROCESO" slot) picks up the visual
# separator style used elsewhere in the column (copy format only, so the
# cell keeps its empty value).
$ws.Range("B9").Copy()
$ws.Range("B38").PasteSpecial(-4122)

# Reflect the new active cell in the bottom-right frozen pane.
$ws.Range("B38").Select()
